$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update player names in column A (rows 2-7) ---
$ws.Range("A2").Value = "Swarabh"
$ws.Range("A3").Value = "Samesh"
$ws.Range("A4").Value = "Aruna"
$ws.Range("A5").Value = "Kathir"
$ws.Range("A6").Value = "Shivam"
$ws.Range("A7").Value = "Abi"

# --- Normalize B14's look to match the rest of column B (hyperlink font + grey fill) ---
# before this, B14 used a style without the grey fill; copy formatting (not value)
# from a sibling cell that already has the right look.
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# --- Every row's email now points at the same address ---
$ws.Range("B2:B14").Value = "aruncyclopse007@gmail.com"

# --- Rebuild the hyperlinks: one on B2, one spanning B3:B14 ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:aruncyclopse007@gmail.com", "", "", "aruncyclopse007@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3:B14"), "mailto:aruncyclopse007@gmail.com", "", "", "aruncyclopse007@gmail.com")

# Hyperlinks.Add() re-stamps the top-left cell of each range with a slightly
# different (but visually identical) style; pull the plain B2:B14 look back
# from an unaffected sibling so every cell in the column matches again.
$ws.Range("B4").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)

# --- Selection moves to the email column ---
$null = $ws.Range("B2:B14").Select()

Write-Output "PlayerDetails updated: names refreshed, emails consolidated, hyperlinks rebuilt."
